$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 2 and row 3: A, B, E, F, G, I, M
# Swap their values between the two rows.
$cols = @("A", "B", "E", "F", "G", "I", "M")

foreach ($col in $cols) {
    $addr2 = "${col}2"
    $addr3 = "${col}3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}
